$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.81973447939941
$ws.Range("C2").Value = 6.851500208941467
$ws.Range("E2").Value = 14.03494175726826
$ws.Range("F2").Value = 42.56847425922016
$ws.Range("G2").Value = 45.42824687196966
$ws.Range("H2").Value = 18.51602644380066
$ws.Range("I2").Value = 29.52213333879394
$ws.Range("J2").Value = 8.959534172934664
$ws.Range("K2").Value = 10.91865707036225
$ws.Range("L2").Value = 12.13910690196044

# Row 3
$ws.Range("B3").Value = 13.61852600825049
$ws.Range("C3").Value = 6.794156067259964
$ws.Range("E3").Value = 14.0324622594094
$ws.Range("F3").Value = 42.57776639116585
$ws.Range("G3").Value = 45.50580613379724
$ws.Range("H3").Value = 18.56781521643094
$ws.Range("I3").Value = 29.61013600299996
$ws.Range("J3").Value = 8.967133985197991
$ws.Range("K3").Value = 10.78465274997188
$ws.Range("L3").Value = 12.12475469688944

# Row 4
$ws.Range("B4").Value = 13.49629444749501
$ws.Range("C4").Value = 6.757987170090307
$ws.Range("E4").Value = 14.03281036958943
$ws.Range("F4").Value = 42.59275964916917
$ws.Range("G4").Value = 45.56504843852787
$ws.Range("H4").Value = 18.60250019668616
$ws.Range("I4").Value = 29.66931850493536
$ws.Range("J4").Value = 8.972044211809447
$ws.Range("K4").Value = 10.70369614789093
$ws.Range("L4").Value = 12.11771419422659

# Row 5
$ws.Range("B5").Value = 13.44687588348699
$ws.Range("C5").Value = 6.743010746399198
$ws.Range("E5").Value = 14.03342422168729
$ws.Range("F5").Value = 42.60120507607537
$ws.Range("G5").Value = 45.59210291821352
$ws.Range("H5").Value = 18.61735991233787
$ws.Range("I5").Value = 29.6947287156595
$ws.Range("J5").Value = 8.974106717869605
$ws.Range("K5").Value = 10.67107586279779
$ws.Range("L5").Value = 12.11529325530984

# Row 6
$ws.Range("B6").Value = 13.43869543874442
$ws.Range("C6").Value = 6.740509706070476
$ws.Range("E6").Value = 14.03355470339135
$ws.Range("F6").Value = 42.60274847289782
$ws.Range("G6").Value = 45.59677092421995
$ws.Range("H6").Value = 18.61987114791632
$ws.Range("I6").Value = 29.69902609467317
$ws.Range("J6").Value = 8.974452919163161
$ws.Range("K6").Value = 10.6656827347965
$ws.Range("L6").Value = 12.11491839669034

# Row 7
$ws.Range("B7").Value = 13.49562630336874
$ws.Range("C7").Value = 6.757786148393641
$ws.Range("E7").Value = 14.03281673507341
$ws.Range("F7").Value = 42.5928640914312
$ws.Range("G7").Value = 45.56540152423215
$ws.Range("H7").Value = 18.60269766399973
$ws.Range("I7").Value = 29.66965596359282
$ws.Range("J7").Value = 8.972071777997332
$ws.Range("K7").Value = 10.70325467121032
$ws.Range("L7").Value = 12.1176797269206

# Row 8
$ws.Range("B8").Value = 13.75012294929635
$ws.Range("C8").Value = 6.831929401680825
$ws.Range("E8").Value = 14.03369962009584
$ws.Range("F8").Value = 42.56975085535301
$ws.Range("G8").Value = 45.45257213997888
$ws.Range("H8").Value = 18.53328381959142
$ws.Range("I8").Value = 29.55140703212744
$ws.Range("J8").Value = 8.962104094246815
$ws.Range("K8").Value = 10.87220060838584
$ws.Range("L8").Value = 12.13379214253872

# Row 9
$ws.Range("B9").Value = 14.25686562088248
$ws.Range("C9").Value = 6.969544756225366
$ws.Range("E9").Value = 14.05019618698364
$ws.Range("F9").Value = 42.59807404266843
$ws.Range("G9").Value = 45.32389246444409
$ws.Range("H9").Value = 18.42008921091102
$ws.Range("I9").Value = 29.36045036788814
$ws.Range("J9").Value = 8.94448301073545
$ws.Range("K9").Value = 11.21232300809525
$ws.Range("L9").Value = 12.17932463663682

# Row 10
$ws.Range("B10").Value = 14.63022633506535
$ws.Range("C10").Value = 7.065692713333474
$ws.Range("E10").Value = 14.07120759383872
$ws.Range("F10").Value = 42.66365168208365
$ws.Range("G10").Value = 45.28623619088423
$ws.Range("H10").Value = 18.35093210900564
$ws.Range("I10").Value = 29.24521072792343
$ws.Range("J10").Value = 8.932696935380269
$ws.Range("K10").Value = 11.46533436110275
$ws.Range("L10").Value = 12.22109634027501

# Row 11
$ws.Range("B11").Value = 14.79955957400067
$ws.Range("C11").Value = 7.108312189367057
$ws.Range("E11").Value = 14.08267027411298
$ws.Range("F11").Value = 42.7031526164087
$ws.Range("G11").Value = 45.28152847264511
$ws.Range("H11").Value = 18.32251858036767
$ws.Range("I11").Value = 29.19824748409272
$ws.Range("J11").Value = 8.927584140963679
$ws.Range("K11").Value = 11.5806392575804
$ws.Range("L11").Value = 12.2418632270663

# Row 12
$ws.Range("B12").Value = 14.86354800949102
$ws.Range("C12").Value = 7.124286504252404
$ws.Range("E12").Value = 14.08728219115262
$ws.Range("F12").Value = 42.71949427152097
$ws.Range("G12").Value = 45.28153573165257
$ws.Range("H12").Value = 18.3121977043101
$ws.Range("I12").Value = 29.18125056620057
$ws.Range("J12").Value = 8.925683606737962
$ws.Range("K12").Value = 11.62429328248188
$ws.Range("L12").Value = 12.24997685435138

# Row 13
$ws.Range("B13").Value = 14.84977388074192
$ws.Range("C13").Value = 7.120853538149374
$ws.Range("E13").Value = 14.08627691182759
$ws.Range("F13").Value = 42.71591339597418
$ws.Range("G13").Value = 45.28145451485958
$ws.Range("H13").Value = 18.3144009678545
$ws.Range("I13").Value = 29.18487612270274
$ws.Range("J13").Value = 8.926091341966909
$ws.Range("K13").Value = 11.61489264569615
$ws.Range("L13").Value = 12.24821840380244

# Row 14
$ws.Range("B14").Value = 14.80482697634517
$ws.Range("C14").Value = 7.109629729738685
$ws.Range("E14").Value = 14.08304427863561
$ws.Range("F14").Value = 42.70446937096102
$ws.Range("G14").Value = 45.28149317998314
$ws.Range("H14").Value = 18.32166067966599
$ws.Range("I14").Value = 29.19683335313162
$ws.Range("J14").Value = 8.927427071132389
$ws.Range("K14").Value = 11.58423111249355
$ws.Range("L14").Value = 12.24252576089829

# Row 15
$ws.Range("B15").Value = 14.77727642373364
$ws.Range("C15").Value = 7.102733253227032
$ws.Range("E15").Value = 14.08109944165287
$ws.Range("F15").Value = 42.69763952095165
$ws.Range("G15").Value = 45.28175005617907
$ws.Range("H15").Value = 18.32616461454752
$ws.Range("I15").Value = 29.2042600542372
$ws.Range("J15").Value = 8.928249870177627
$ws.Range("K15").Value = 11.56544765293163
$ws.Range("L15").Value = 12.23907123972311

# Row 16
$ws.Range("B16").Value = 14.61914447871018
$ws.Range("C16").Value = 7.062884584645909
$ws.Range("E16").Value = 14.07049658420462
$ws.Range("F16").Value = 42.66126421223835
$ws.Range("G16").Value = 45.28679416287001
$ws.Range("H16").Value = 18.35285035657306
$ws.Range("I16").Value = 29.24838991843284
$ws.Range("J16").Value = 8.933036056885978
$ws.Range("K16").Value = 11.45779972367917
$ws.Range("L16").Value = 12.21977434039273

# Row 17
$ws.Range("B17").Value = 14.52196089188174
$ws.Range("C17").Value = 7.038149393538684
$ws.Range("E17").Value = 14.06447805220419
$ws.Range("F17").Value = 42.64142117075742
$ws.Range("G17").Value = 45.29307323274339
$ws.Range("H17").Value = 18.37000191288273
$ws.Range("I17").Value = 29.27686194865836
$ws.Range("J17").Value = 8.936035795601349
$ws.Range("K17").Value = 11.39178634148828
$ws.Range("L17").Value = 12.20838545839692

# Row 18
$ws.Range("B18").Value = 14.46601894903898
$ws.Range("C18").Value = 7.023817294113939
$ws.Range("E18").Value = 14.06119577803499
$ws.Range("F18").Value = 42.63091859936189
$ws.Range("G18").Value = 45.29785389863301
$ws.Range("H18").Value = 18.38015373734328
$ws.Range("I18").Value = 29.29375210292902
$ws.Range("J18").Value = 8.937784591371036
$ws.Range("K18").Value = 11.35383915390678
$ws.Range("L18").Value = 12.2020012401405

# Row 19
$ws.Range("B19").Value = 14.44707221444839
$ws.Range("C19").Value = 7.018946765269265
$ws.Range("E19").Value = 14.060115349941
$ws.Range("F19").Value = 42.62751919908624
$ws.Range("G19").Value = 45.29967321284863
$ws.Range("H19").Value = 18.38364019218672
$ws.Range("I19").Value = 29.29955900415402
$ws.Range("J19").Value = 8.93838073302866
$ws.Range("K19").Value = 11.34099587637986
$ws.Range("L19").Value = 12.19986834831314

# Row 20
$ws.Range("B20").Value = 14.53231130592726
$ws.Range("C20").Value = 7.040793399848226
$ws.Range("E20").Value = 14.06510018368081
$ws.Range("F20").Value = 42.64343929654854
$ws.Range("G20").Value = 45.29228378566602
$ws.Range("H20").Value = 18.36814642299136
$ws.Range("I20").Value = 29.27377786233775
$ws.Range("J20").Value = 8.935714045301504
$ws.Range("K20").Value = 11.39881159404815
$ws.Range("L20").Value = 12.20958063515303

# Row 21
$ws.Range("B21").Value = 14.81803309463847
$ws.Range("C21").Value = 7.112930936261687
$ws.Range("E21").Value = 14.08398644140537
$ws.Range("F21").Value = 42.70779327217026
$ws.Range("G21").Value = 45.28143322129513
$ws.Range("H21").Value = 18.31951641483289
$ws.Range("I21").Value = 29.19329984957981
$ws.Range("J21").Value = 8.92703377132255
$ws.Range("K21").Value = 11.59323771425225
$ws.Range("L21").Value = 12.24419108634803

# Row 22
$ws.Range("B22").Value = 15.00395872409679
$ws.Range("C22").Value = 7.15911492952083
$ws.Range("E22").Value = 14.09790941885148
$ws.Range("F22").Value = 42.75791233041177
$ws.Range("G22").Value = 45.28477581378901
$ws.Range("H22").Value = 18.29029131598346
$ws.Range("I22").Value = 29.14529115046694
$ws.Range("J22").Value = 8.921567948303826
$ws.Range("K22").Value = 11.72023377441398
$ws.Range("L22").Value = 12.26826414787342

# Row 23
$ws.Range("B23").Value = 14.90482043694721
$ws.Range("C23").Value = 7.134554950460073
$ws.Range("E23").Value = 14.09033482595378
$ws.Range("F23").Value = 42.73042797932624
$ws.Range("G23").Value = 45.28203622933411
$ws.Range("H23").Value = 18.30565509035845
$ws.Range("I23").Value = 29.1704938459629
$ws.Range("J23").Value = 8.92446626431574
$ws.Range("K23").Value = 11.65247303391345
$ws.Range("L23").Value = 12.25528435664999

# Row 24
$ws.Range("B24").Value = 14.52763209648478
$ws.Range("C24").Value = 7.039598391660198
$ws.Range("E24").Value = 14.06481836384289
$ws.Range("F24").Value = 42.64252408088937
$ws.Range("G24").Value = 45.29263704809618
$ws.Range("H24").Value = 18.36898438301623
$ws.Range("I24").Value = 29.2751705543332
$ws.Range("J24").Value = 8.935859433187659
$ws.Range("K24").Value = 11.39563545943213
$ws.Range("L24").Value = 12.20903978646037

# Row 25
$ws.Range("B25").Value = 14.11933399943709
$ws.Range("C25").Value = 6.93317151125405
$ws.Range("E25").Value = 14.04416339971036
$ws.Range("F25").Value = 42.58253569198573
$ws.Range("G25").Value = 45.34874179312601
$ws.Range("H25").Value = 18.44825373873537
$ws.Range("I25").Value = 29.40771643616787
$ws.Range("J25").Value = 8.949045264198638
$ws.Range("K25").Value = 11.11959965379881
$ws.Range("L25").Value = 12.16553222116168

